$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3284.389
$ws.Range("I64").Value = 3148
$ws.Range("J64").Value = 3761.75
$ws.Range("K64").Value = 3148
$ws.Range("L64").Value = 3761.75
$ws.Range("M64").Value = -2900
$ws.Range("N64").Value = -4257.75
$ws.Range("H67").Value = 3284.389
$ws.Range("I67").Value = 3148
$ws.Range("J67").Value = 3761.75
$ws.Range("K67").Value = 3148
$ws.Range("L67").Value = 3761.75
$ws.Range("M67").Value = -2290
$ws.Range("N67").Value = -5477.75
$ws.Range("H137").Value = 1037.9697
$ws.Range("I137").Value = 921.9
$ws.Range("J137").Value = 1216.5385
$ws.Range("K137").Value = 2765.7
$ws.Range("L137").Value = 3649.6155
$ws.Range("M137").Value = -215.6999999999998
$ws.Range("N137").Value = -8749.6155
$ws.Range("H139").Value = 70180
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 70180
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 70180
$ws.Range("N139").Value = -80460
$ws.Range("H140").Value = 88250
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 88250
$ws.Range("K140").Value = 0
$ws.Range("L140").ClearContents()
$ws.Range("M140").Value = 88250
$ws.Range("N140").Value = -98610

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3481.225
$ws.Range("I61").Value = 3675.303
$ws.Range("J61").Value = 2566.2856
$ws.Range("K61").Value = 3675.303
$ws.Range("L61").Value = 2566.2856
$ws.Range("M61").Value = -3463.303
$ws.Range("N61").Value = -2990.2856
$ws.Range("H74").Value = 1413.5
$ws.Range("I74").Value = 843
$ws.Range("J74").Value = 3125
$ws.Range("K74").Value = 843
$ws.Range("L74").Value = 3125
$ws.Range("M74").Value = 31
$ws.Range("N74").Value = -4873
$ws.Range("H77").Value = 1413.5
$ws.Range("I77").Value = 843
$ws.Range("J77").Value = 3125
$ws.Range("K77").Value = 4215
$ws.Range("L77").Value = 15625
$ws.Range("M77").Value = 153
$ws.Range("N77").Value = -24361
$ws.Range("H136").Value = 3481.225
$ws.Range("I136").Value = 3675.303
$ws.Range("J136").Value = 2566.2856
$ws.Range("K136").Value = 11025.909
$ws.Range("L136").Value = 7698.8568
$ws.Range("M136").Value = -8475.909
$ws.Range("N136").Value = -12798.8568
$ws.Range("H138").Value = 62660
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 62660
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 62660
$ws.Range("N138").Value = -72940
$ws.Range("H139").Value = 53238.332
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 53238.332
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 53238.332
$ws.Range("N139").Value = -63518.332

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1677.5
$ws.Range("I105").Value = 1675
$ws.Range("J105").Value = 1700
$ws.Range("K105").Value = 1675
$ws.Range("L105").Value = 1700
$ws.Range("M105").Value = 72
$ws.Range("N105").Value = -5194
$ws.Range("H138").Value = 59575
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 59575
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 59575
$ws.Range("N138").Value = -69855
$ws.Range("H140").Value = 89740
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 89740
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 89740
$ws.Range("N140").Value = -100100

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4712.393
$ws.Range("I31").Value = 5224.5264
$ws.Range("J31").Value = 4449.4053
$ws.Range("K31").Value = 5224.5264
$ws.Range("L31").Value = 4449.4053
$ws.Range("M31").Value = -4929.5264
$ws.Range("N31").Value = -5039.4053
$ws.Range("H34").Value = 4712.393
$ws.Range("I34").Value = 5224.5264
$ws.Range("J34").Value = 4449.4053
$ws.Range("K34").Value = 5224.5264
$ws.Range("L34").Value = 4449.4053
$ws.Range("M34").Value = -5022.5264
$ws.Range("N34").Value = -4853.4053
$ws.Range("H62").Value = 3146.7273
$ws.Range("I62").Value = 3169.3333
$ws.Range("J62").Value = 3045
$ws.Range("K62").Value = 3169.3333
$ws.Range("L62").Value = 3045
$ws.Range("M62").Value = -2545.3333
$ws.Range("N62").Value = -4293
$ws.Range("H65").Value = 3146.7273
$ws.Range("I65").Value = 3169.3333
$ws.Range("J65").Value = 3045
$ws.Range("K65").Value = 15846.6665
$ws.Range("L65").Value = 15225
$ws.Range("M65").Value = -12726.6665
$ws.Range("N65").Value = -21465
$ws.Range("H140").Value = 90000
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 90000
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 90000
$ws.Range("N140").Value = -100360

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 100433.336
$ws.Range("I9").Value = 300000
$ws.Range("J9").Value = 650
$ws.Range("K9").Value = 900000
$ws.Range("L9").Value = 1950
$ws.Range("M9").Value = -899776
$ws.Range("N9").Value = -2398
$ws.Range("H22").Value = 29412700
$ws.Range("I22").Value = 979.6
$ws.Range("J22").Value = 41667584
$ws.Range("K22").Value = 2938.8
$ws.Range("L22").Value = 125002752
$ws.Range("M22").Value = -2769.8
$ws.Range("N22").Value = -125003090
$ws.Range("H27").Value = 29412700
$ws.Range("I27").Value = 979.6
$ws.Range("J27").Value = 41667584
$ws.Range("K27").Value = 2938.8
$ws.Range("L27").Value = 125002752
$ws.Range("M27").Value = -2836.8
$ws.Range("N27").Value = -125002956
$ws.Range("H40").Value = 257.7742
$ws.Range("I40").Value = 145.96297
$ws.Range("J40").Value = 1012.5
$ws.Range("K40").Value = 583.8518800000001
$ws.Range("L40").Value = 4050
$ws.Range("M40").Value = -514.8518800000001
$ws.Range("N40").Value = -4188
$ws.Range("H131").Value = 7813465.5
$ws.Range("I131").Value = 1437
$ws.Range("J131").Value = 9260137
$ws.Range("K131").Value = 4311
$ws.Range("L131").Value = 27780411
$ws.Range("M131").Value = 729
$ws.Range("N131").Value = -27790491
$ws.Range("H133").Value = 7685.3267
$ws.Range("I133").Value = 5997.615
$ws.Range("J133").Value = 8294.777
$ws.Range("K133").Value = 17992.845
$ws.Range("L133").Value = 24884.331
$ws.Range("M133").Value = -12932.845
$ws.Range("N133").Value = -35004.331
$ws.Range("H139").Value = 2412
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 2412
$ws.Range("K139").Value = 0
$ws.Range("L139").ClearContents()
$ws.Range("M139").Value = 7236
$ws.Range("N139").Value = -17516

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14430277
$ws.Range("I70").Value = 20839538
$ws.Range("J70").Value = 9443.166999999999
$ws.Range("K70").Value = 20839538
$ws.Range("L70").Value = 9443.166999999999
$ws.Range("M70").Value = -20839268
$ws.Range("N70").Value = -9983.166999999999
$ws.Range("H73").Value = 14430277
$ws.Range("I73").Value = 20839538
$ws.Range("J73").Value = 9443.166999999999
$ws.Range("K73").Value = 20839538
$ws.Range("L73").Value = 9443.166999999999
$ws.Range("M73").Value = -20838602
$ws.Range("N73").Value = -11315.167
$ws.Range("H80").Value = 2315
$ws.Range("I80").Value = 2341
$ws.Range("J80").Value = 2250
$ws.Range("K80").Value = 2341
$ws.Range("L80").Value = 2250
$ws.Range("M80").Value = -1343
$ws.Range("N80").Value = -4246
$ws.Range("H83").Value = 2315
$ws.Range("I83").Value = 2341
$ws.Range("J83").Value = 2250
$ws.Range("K83").Value = 11705
$ws.Range("L83").Value = 11250
$ws.Range("M83").Value = -6713
$ws.Range("N83").Value = -21234
$ws.Range("H138").Value = 68066.664
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 68066.664
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 68066.664
$ws.Range("N138").Value = -78346.664
$ws.Range("H140").Value = 99893
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 99893
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 99893
$ws.Range("N140").Value = -110253
$ws.Range("H141").Value = 67990
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 67990
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 67990
$ws.Range("N141").Value = -78350

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5498.278
$ws.Range("I136").Value = 3348.5454
$ws.Range("J136").Value = 8876.429
$ws.Range("K136").Value = 10045.6362
$ws.Range("L136").Value = 26629.287
$ws.Range("M136").Value = -7495.636200000001
$ws.Range("N136").Value = -31729.287

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1020.29034
$ws.Range("I136").Value = 972.087
$ws.Range("J136").Value = 1158.875
$ws.Range("K136").Value = 2916.261
$ws.Range("L136").Value = 3476.625
$ws.Range("M136").Value = -366.261
$ws.Range("N136").Value = -8576.625
$ws.Range("H138").Value = 62033.332
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 62033.332
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 62033.332
$ws.Range("N138").Value = -72313.33199999999
$ws.Range("H139").Value = 54683.332
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 54683.332
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 54683.332
$ws.Range("N139").Value = -64963.332
